$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking Price cells to stay as text (matches source data which
# stores these as literal strings, e.g. "538.33" rather than numeric 538.33).
foreach ($cell in @("D5", "D6", "D14", "D19", "D21", "D24", "D26", "D27", "D30", "D32", "D35", "D41", "D42", "D47", "D49")) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '58.696.24'
$ws.Range("D3").Value = '2.309.65'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '538.33'
$ws.Range("E5").Value = '  -2.19%  '
$ws.Range("D6").Value = '132.35'
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D9").Value = '2.308.49'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").Value = '23.80'
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("D15").Value = '2.720.75'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").Value = '58.593.39'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '2.318.19'
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").Value = '10.61'
$ws.Range("E19").Value = '  -1.04%  '
$ws.Range("E20").Value = '  -3.11%  '
$ws.Range("D21").Value = '315.82'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '63.22'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = '7.94'
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("E29").Value = '  -2.06%  '
$ws.Range("D30").Value = '171.18'
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").Value = '1.10'
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("E34").Value = '  +0.39%  '
$ws.Range("D35").Value = '17.95'
$ws.Range("E35").Value = '  +0.60%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").Value = '291.39'
$ws.Range("E41").Value = '  -5.02%  '
$ws.Range("D42").Value = '141.21'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").Value = '18.36'
$ws.Range("E47").Value = '  -2.62%  '
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").Value = '10.95'
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("E51").Value = '  +0.01%  '

# Restore default (Normal) style on those cells so only the value - not formatting -
# differs from the original workbook.
foreach ($cell in @("D5", "D6", "D14", "D19", "D21", "D24", "D26", "D27", "D30", "D32", "D35", "D41", "D42", "D47", "D49")) {
    $ws.Range($cell).Style = "Normal"
}
